# Update cryptos list: refresh prices/volume %, and insert new "RenzoRestakedETH" row
# at position 37 (pushing Hedera..Stacks down by one row and dropping EnergySwap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (many of these look like numbers, e.g. "1.00"
# or "0.0000154", and Excel would silently coerce them to numeric cells and lose
# the original formatting). Forcing the cell to the "@" (Text) number format first
# keeps them as text; the original style is restored afterwards so no stray
# formatting is left behind.
function Set-TextValue($Worksheet, $Row, $Col, $Text) {
    $cell = $Worksheet.Cells.Item($Row, $Col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws 2 4 '67.711.75'
Set-TextValue $ws 2 5 '  -0.24%  '

# Row 3
Set-TextValue $ws 3 4 '3.800.88'
Set-TextValue $ws 3 5 '  +0.65%  '

# Row 4
Set-TextValue $ws 4 5 '  -0.14%  '

# Row 5
Set-TextValue $ws 5 4 '596.39'
Set-TextValue $ws 5 5 '  +0.50%  '

# Row 6
Set-TextValue $ws 6 4 '167.17'
Set-TextValue $ws 6 5 '  +0.37%  '

# Row 7
Set-TextValue $ws 7 4 '3.798.33'
Set-TextValue $ws 7 5 '  +0.59%  '

# Row 8
Set-TextValue $ws 8 5 '  +0.02%  '

# Row 9
Set-TextValue $ws 9 5 '  +0.42%  '

# Row 10
Set-TextValue $ws 10 5 '  +1.05%  '

# Row 11
Set-TextValue $ws 11 5 '  -1.19%  '

# Row 12
Set-TextValue $ws 12 5 '  -0.27%  '

# Row 13
Set-TextValue $ws 13 5 '  -1.62%  '

# Row 14
Set-TextValue $ws 14 4 '36.09'
Set-TextValue $ws 14 5 '  +0.34%  '

# Row 15
Set-TextValue $ws 15 4 '4.439.40'
Set-TextValue $ws 15 5 '  +0.41%  '

# Row 16
Set-TextValue $ws 16 4 '3.802.92'
Set-TextValue $ws 16 5 '  +0.19%  '

# Row 17
Set-TextValue $ws 17 4 '18.58'
Set-TextValue $ws 17 5 '  +4.08%  '

# Row 18
Set-TextValue $ws 18 4 '67.726.18'
Set-TextValue $ws 18 5 '  -0.42%  '

# Row 19
Set-TextValue $ws 19 4 '7.11'
Set-TextValue $ws 19 5 '  +2.43%  '

# Row 21
Set-TextValue $ws 21 4 '461.26'
Set-TextValue $ws 21 5 '  +0.02%  '

# Row 22
Set-TextValue $ws 22 4 '9.93'
Set-TextValue $ws 22 5 '  -5.62%  '

# Row 23
Set-TextValue $ws 23 5 '  +0.95%  '

# Row 24
Set-TextValue $ws 24 4 '0.0000154'
Set-TextValue $ws 24 5 '  +1.32%  '

# Row 25
Set-TextValue $ws 25 4 '83.52'
Set-TextValue $ws 25 5 '  +0.07%  '

# Row 26
Set-TextValue $ws 26 5 '  +2.12%  '

# Row 27
Set-TextValue $ws 27 4 '2.11'
Set-TextValue $ws 27 5 '  -1.68%  '

# Row 28
Set-TextValue $ws 28 4 '1.01'
Set-TextValue $ws 28 5 '  +0.68%  '

# Row 29
Set-TextValue $ws 29 4 '10.00'
Set-TextValue $ws 29 5 '  -0.27%  '

# Row 30
Set-TextValue $ws 30 4 '3.949.88'
Set-TextValue $ws 30 5 '  +0.50%  '

# Row 31
Set-TextValue $ws 31 5 '  +1.43%  '

# Row 32
Set-TextValue $ws 32 4 '2.26'
Set-TextValue $ws 32 5 '  +5.21%  '

# Row 33
Set-TextValue $ws 33 4 '7.24'
Set-TextValue $ws 33 5 '  +0.79%  '

# Row 34
Set-TextValue $ws 34 4 '29.72'
Set-TextValue $ws 34 5 '  -0.31%  '

# Row 35
Set-TextValue $ws 35 5 '  +0.02%  '

# Row 36
Set-TextValue $ws 36 5 '  +0.32%  '

# Row 37
Set-TextValue $ws 37 2 'RenzoRestakedETH'
Set-TextValue $ws 37 3 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue $ws 37 4 '3.742.64'
Set-TextValue $ws 37 5 '  +0.35%  '

# Row 38
Set-TextValue $ws 38 2 'Hedera'
Set-TextValue $ws 38 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 38 4 '0.100'
Set-TextValue $ws 38 5 '  +0.03%  '

# Row 39
Set-TextValue $ws 39 2 'dogwifhat'
Set-TextValue $ws 39 3 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 39 4 '3.38'
Set-TextValue $ws 39 5 '  -1.40%  '

# Row 40
Set-TextValue $ws 40 2 'Kaspa'
Set-TextValue $ws 40 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 40 4 '0.138'
Set-TextValue $ws 40 5 '  +0.28%  '

# Row 41
Set-TextValue $ws 41 2 'Mantle'
Set-TextValue $ws 41 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws 41 4 '0.995'
Set-TextValue $ws 41 5 '  +0.14%  '

# Row 42
Set-TextValue $ws 42 2 'Filecoin'
Set-TextValue $ws 42 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 42 4 '5.79'
Set-TextValue $ws 42 5 '  +0.93%  '

# Row 43
Set-TextValue $ws 43 2 'FirstDigitalUSD'
Set-TextValue $ws 43 3 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 43 4 '0.999'
Set-TextValue $ws 43 5 '  -0.13%  '

# Row 44
Set-TextValue $ws 44 2 'USDe'
Set-TextValue $ws 44 3 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws 44 4 '1.00'
Set-TextValue $ws 44 5 '  +0.01%  '

# Row 45
Set-TextValue $ws 45 2 'OKB'
Set-TextValue $ws 45 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 45 4 '48.09'
Set-TextValue $ws 45 5 '  +2.71%  '

# Row 46
Set-TextValue $ws 46 2 'Arweave'
Set-TextValue $ws 46 3 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws 46 4 '43.86'
Set-TextValue $ws 46 5 '  +1.52%  '

# Row 47
Set-TextValue $ws 47 2 'TheGraph'
Set-TextValue $ws 47 3 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws 47 4 '0.298'
Set-TextValue $ws 47 5 '  -0.15%  '

# Row 48
Set-TextValue $ws 48 2 'Monero'
Set-TextValue $ws 48 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 48 4 '150.07'
Set-TextValue $ws 48 5 '  +2.15%  '

# Row 49
Set-TextValue $ws 49 2 'Cosmos'
Set-TextValue $ws 49 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws 49 4 '8.33'
Set-TextValue $ws 49 5 '  -0.11%  '

# Row 50
Set-TextValue $ws 50 2 'Bittensor'
Set-TextValue $ws 50 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws 50 4 '393.97'
Set-TextValue $ws 50 5 '  +1.91%  '

# Row 51
Set-TextValue $ws 51 2 'Stacks'
Set-TextValue $ws 51 3 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws 51 4 '1.83'
Set-TextValue $ws 51 5 '  -3.95%  '
